# week 3 data updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 25 (Prism, Week 3): fill in the previously-missing SF1..SF36 / PSS1..PSS10 values ---
# Columns C (3) through AV (48); B25 (Week=3) already has its value.
$row25Values = @(2,1,3,3,3,3,3,3,3,3,3,3,2,2,1,2,2,2,2,1,1,1,1,3,6,2,2,6,4,1,4,5,3,2,5,1,1,0,2,4,4,1,4,4,0,0)

for ($i = 0; $i -lt $row25Values.Length; $i++) {
    $col = 3 + $i
    $ws.Cells.Item(25, $col).Value = $row25Values[$i]
}

# --- Row 31 (Meadow, Week 3): correct AR31, and add AU31 / AV31 ---
$ws.Range("AR31").Value = 1
$ws.Range("AU31").Value = 1
$ws.Range("AV31").Value = 0

# --- Update the active selection to match the saved view (A31, no horizontal scroll) ---
$ws.Range("A31").Select() | Out-Null
